$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$p = $d.Paragraphs.Item(1)
$pr = $p.Range
$pr.MoveEnd(1, -1) | Out-Null
$pr.Text = "2023-06-29 Thursday"

# Update each multiplication-problem cell in the table, addressed by row/column
# (some old values repeat across cells, so each cell is targeted individually
#  via its own Range rather than a document-wide Find/Replace)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "11×62="
$cell = $t.Cell(1, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "41×28="
$cell = $t.Cell(1, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "51×46="
$cell = $t.Cell(1, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "28×13="
$cell = $t.Cell(1, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "25×21="

$cell = $t.Cell(2, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "59×74="
$cell = $t.Cell(2, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "47×34="
$cell = $t.Cell(2, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "19×97="
$cell = $t.Cell(2, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "60×93="
$cell = $t.Cell(2, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "99×83="

$cell = $t.Cell(3, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "50×19="
$cell = $t.Cell(3, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "43×22="
$cell = $t.Cell(3, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "17×37="
$cell = $t.Cell(3, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "15×13="
$cell = $t.Cell(3, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "60×26="

$cell = $t.Cell(4, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "79×94="
$cell = $t.Cell(4, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "63×85="
$cell = $t.Cell(4, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "93×50="
$cell = $t.Cell(4, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "86×21="
$cell = $t.Cell(4, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "39×63="

$cell = $t.Cell(5, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "36×31="
$cell = $t.Cell(5, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "22×15="
$cell = $t.Cell(5, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "69×100="
$cell = $t.Cell(5, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "21×23="
$cell = $t.Cell(5, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "33×90="

$cell = $t.Cell(6, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "21×82="
$cell = $t.Cell(6, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "38×70="
$cell = $t.Cell(6, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "43×78="
$cell = $t.Cell(6, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "32×87="
$cell = $t.Cell(6, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "55×87="

$cell = $t.Cell(7, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "97×78="
$cell = $t.Cell(7, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "73×59="
$cell = $t.Cell(7, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "74×42="
$cell = $t.Cell(7, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "94×24="
$cell = $t.Cell(7, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "27×30="

$cell = $t.Cell(8, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "72×77="
$cell = $t.Cell(8, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "87×81="
$cell = $t.Cell(8, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "27×68="
$cell = $t.Cell(8, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "95×79="
$cell = $t.Cell(8, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "44×18="

$cell = $t.Cell(9, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "98×100="
$cell = $t.Cell(9, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "53×99="
$cell = $t.Cell(9, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "46×77="
$cell = $t.Cell(9, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "20×61="
$cell = $t.Cell(9, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "58×13="

$cell = $t.Cell(10, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "54×82="
$cell = $t.Cell(10, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "56×29="
$cell = $t.Cell(10, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "60×19="
$cell = $t.Cell(10, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "18×54="
$cell = $t.Cell(10, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "79×45="

$cell = $t.Cell(11, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "94×69="
$cell = $t.Cell(11, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "19×64="
$cell = $t.Cell(11, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "46×68="
$cell = $t.Cell(11, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "30×21="
$cell = $t.Cell(11, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "71×61="

$cell = $t.Cell(12, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "68×22="
$cell = $t.Cell(12, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "77×39="
$cell = $t.Cell(12, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "100×43="
$cell = $t.Cell(12, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "96×32="
$cell = $t.Cell(12, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "74×69="

$cell = $t.Cell(13, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "46×63="
$cell = $t.Cell(13, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "89×74="
$cell = $t.Cell(13, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "67×55="
$cell = $t.Cell(13, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "41×48="
$cell = $t.Cell(13, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "48×90="

$cell = $t.Cell(14, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "93×78="
$cell = $t.Cell(14, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "95×91="
$cell = $t.Cell(14, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "22×52="
$cell = $t.Cell(14, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "64×18="
$cell = $t.Cell(14, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "99×80="

$cell = $t.Cell(15, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "19×40="
$cell = $t.Cell(15, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "65×65="
$cell = $t.Cell(15, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "52×40="
$cell = $t.Cell(15, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "29×58="
$cell = $t.Cell(15, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "75×58="

$cell = $t.Cell(16, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "26×48="
$cell = $t.Cell(16, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "60×79="
$cell = $t.Cell(16, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "65×82="
$cell = $t.Cell(16, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "65×76="
$cell = $t.Cell(16, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "11×73="

$cell = $t.Cell(17, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "87×26="
$cell = $t.Cell(17, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "44×54="
$cell = $t.Cell(17, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "24×19="
$cell = $t.Cell(17, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "31×54="
$cell = $t.Cell(17, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "54×95="

$cell = $t.Cell(18, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "34×17="
$cell = $t.Cell(18, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "61×93="
$cell = $t.Cell(18, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "20×16="
$cell = $t.Cell(18, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "15×51="
$cell = $t.Cell(18, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "17×86="

$cell = $t.Cell(19, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "33×14="
$cell = $t.Cell(19, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "46×44="
$cell = $t.Cell(19, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "93×96="
$cell = $t.Cell(19, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "37×64="
$cell = $t.Cell(19, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "42×59="

$cell = $t.Cell(20, 1)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "62×96="
$cell = $t.Cell(20, 2)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "26×37="
$cell = $t.Cell(20, 3)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "16×74="
$cell = $t.Cell(20, 4)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "55×97="
$cell = $t.Cell(20, 5)
$cr = $cell.Range
$cr.MoveEnd(1, -1) | Out-Null
$cr.Text = "31×35="

